# Generate Report for Handoff
# The localization file "f86032e9-c9e9-421d-a86d-6d774bc12014.md" moved from
# "Handed back: in sync with en-US" to "Ready for handoff" status, with the
# corresponding handoff timestamps bumped on the Overview / zh-cn / de-de
# sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 is the f86032e9-...-.md file ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"
$overview.Range("D3").Value = "2016-32-12 00:32:56"

# --- zh-cn sheet: row 3 is the f86032e9-...-.md file ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("E3").Value = "2016-03-12 00:32:53"

# --- de-de sheet: row 3 is the f86032e9-...-.md file ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("E3").Value = "2016-03-12 00:32:56"
